$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update department (C2) from "FACULTY OF ENGLISH" to "English"
$ws.Range("C2").Value = "English"

# Clear promotionValidity (R2) value
$ws.Range("R2").Value = ""

# Row height shrinks since R2's long wrapped text no longer needs as much
# vertical space once the promotionValidity text is removed
$ws.Rows.Item(2).RowHeight = 28.5
